$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.384.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.626.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.16%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.035.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.625.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.338.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0946"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.50%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0807"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "129.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.02%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0332"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.40%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.059.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +8.14%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.64%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.05%  "
$ws.Range("E51").Style = "Normal"
